$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.15188
$ws.Range("H2").Value = 0.45564
$ws.Range("I2").Value = 0.229582042173683
$ws.Range("J2").Value = 0.229582042173683
$ws.Range("M2").Value = 9.101794333333332
$ws.Range("N2").Value = 27.305383
$ws.Range("O2").Value = 0.1526015110517656
$ws.Range("P2").Value = 0.1526015110517656
$ws.Range("Q2").Value = 1.382380523346666
$ws.Range("R2").Value = 12.44142471012
$ws.Range("S2").Value = 0.0350345665460542
$ws.Range("T2").Value = 0.0350345665460542

# Row 3
$ws.Range("G3").Value = 0.15188
$ws.Range("H3").Value = 0.45564
$ws.Range("I3").Value = 0.229582042173683
$ws.Range("J3").Value = 0.229582042173683
$ws.Range("O3").Value = 0.5991759712230392
$ws.Range("P3").Value = 0.5991759712230392
$ws.Range("Q3").Value = 5.427791553093333
$ws.Range("R3").Value = 48.85012397784
$ws.Range("S3").Value = 0.1375600430947853
$ws.Range("T3").Value = 0.1375600430947853

# Row 4
$ws.Range("G4").Value = 0.15188
$ws.Range("H4").Value = 0.45564
$ws.Range("I4").Value = 0.229582042173683
$ws.Range("J4").Value = 0.229582042173683
$ws.Range("O4").Value = 0.2482225177251951
$ws.Range("P4").Value = 0.2482225177251951
$ws.Range("Q4").Value = 2.248588310786667
$ws.Range("R4").Value = 20.23729479708
$ws.Range("S4").Value = 0.05698743253284352
$ws.Range("T4").Value = 0.05698743253284352

# Row 5
$ws.Range("I5").Value = 0.7704179578263169
$ws.Range("J5").Value = 0.7704179578263169
$ws.Range("M5").Value = 9.101794333333332
$ws.Range("N5").Value = 27.305383
$ws.Range("O5").Value = 0.1526015110517656
$ws.Range("P5").Value = 0.1526015110517656
$ws.Range("Q5").Value = 4.63891151787
$ws.Range("R5").Value = 41.75020366083
$ws.Range("S5").Value = 0.1175669445057114
$ws.Range("T5").Value = 0.1175669445057114

# Row 6
$ws.Range("I6").Value = 0.7704179578263169
$ws.Range("J6").Value = 0.7704179578263169
$ws.Range("O6").Value = 0.5991759712230392
$ws.Range("P6").Value = 0.5991759712230392
$ws.Range("S6").Value = 0.4616159281282539
$ws.Range("T6").Value = 0.4616159281282539

# Row 7
$ws.Range("I7").Value = 0.7704179578263169
$ws.Range("J7").Value = 0.7704179578263169
$ws.Range("O7").Value = 0.2482225177251951
$ws.Range("P7").Value = 0.2482225177251951
$ws.Range("S7").Value = 0.1912350851923516
$ws.Range("T7").Value = 0.1912350851923516
